$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance")

$ws.Range("A5").Value = "15/02/2026 03:41:32"
$ws.Range("B5").Value = 60.08
$ws.Range("C5").Value = 2.4
$ws.Range("D5").Value = 57.67
$ws.Range("E5").Value = 0.08
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "GRID"

$co1 = $ws.ChartObjects(1)
$chart1 = $co1.Chart
$ser1 = $chart1.SeriesCollection(1)
$ser1.Formula = '=SERIES("Evolução da Banca",Performance!$A$2:$A$5,Performance!$B$2:$B$5,1)'
$excel.CalculateFull()
$excel.CalculateFullRebuild()
Write-Output "done"
